$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = "./images_eeg/Sphere_CW-3.75_BG-grey_stim-white.png"
    3  = "./images_eeg/Sphere_CCW-3.75_BG-grey_stim-white.png"
    4  = "./images_eeg/Sphere_CW-3.75_BG-grey_stim-white.png"
    5  = "./images_eeg/Sphere_CW-3.75_BG-grey_stim-white.png"
    6  = "./images_eeg/Sphere_Ref_BG-grey_stim-yellow.png"
    7  = "./images_eeg/Sphere_CW-3.75_BG-grey_stim-white.png"
    8  = "./images_eeg/Sphere_CCW-3.75_BG-grey_stim-white.png"
    9  = "./images_eeg/Sphere_CW-3.75_BG-grey_stim-yellow.png"
    10 = "./images_eeg/Sphere_CCW-3.75_BG-grey_stim-white.png"
    11 = "./images_eeg/Sphere_CCW-3.75_BG-grey_stim-white.png"
    12 = "./images_eeg/Sphere_CCW-3.75_BG-grey_stim-yellow.png"
    13 = "./images_eeg/Sphere_Ref_BG-grey_stim-white.png"
    14 = "./images_eeg/Sphere_Ref_BG-grey_stim-white.png"
    15 = "./images_eeg/Sphere_CW-3.75_BG-grey_stim-white.png"
    16 = "./images_eeg/Sphere_Ref_BG-grey_stim-yellow.png"
    17 = "./images_eeg/Sphere_CCW-3.75_BG-grey_stim-white.png"
    18 = "./images_eeg/Sphere_CCW-3.75_BG-grey_stim-white.png"
    19 = "./images_eeg/Sphere_CW-3.75_BG-grey_stim-white.png"
}

foreach ($row in $values.Keys) {
    $ws.Range("A$row").Value = $values[$row]
}
